# Updated cryptos list on Sat Nov 23 13:32:10 UTC 2024 with GitHub Actions
# Refresh Price (D) / Volume(1h) (E) columns with latest scraped values, and
# fix a few rows whose rank swapped with a neighbour (name/link/price/volume
# move together). Price cells that look like plain numbers ("1.55", "2.20",
# "0.0000220", ...) are forced to Text via NumberFormat "@" before the
# assignment so Excel doesn't silently reinterpret/round them as numbers;
# the style is then reset to "Normal" so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.834.31"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "3.415.08"
$ws.Range("E3").Value = "  +3.76%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "258.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "669.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.55"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +10.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.471"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +19.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.09"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +22.88%  "
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("D11").Value = "3.412.69"
$ws.Range("E11").Value = "  +3.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.219"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +13.47%  "
$ws.Range("E14").Value = "  +12.84%  "
$ws.Range("B15").Value = "Toncoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.94%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "98.552.42"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("D17").Value = "4.053.32"
$ws.Range("E17").Value = "  +2.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +34.78%  "
$ws.Range("D19").Value = "3.417.61"
$ws.Range("E19").Value = "  +3.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +15.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "535.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +12.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +14.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000220"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.442"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +53.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +15.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "102.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +16.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.10%  "
$ws.Range("D29").Value = "3.594.44"
$ws.Range("E29").Value = "  +3.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.153"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +16.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +20.50%  "
$ws.Range("B32").Value = "Cronos"
$ws.Range("C32").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.197"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.30%  "
$ws.Range("B33").Value = "Dai"
$ws.Range("C33").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "30.36"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.566"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +24.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +14.52%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +11.57%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.162"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "531.33"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.39"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "24.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0440"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +36.62%  "
$ws.Range("B44").Value = "MantraDAO"
$ws.Range("C44").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.21%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.857"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.95%  "
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +18.76%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +12.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +15.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +17.64%  "
